$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header metadata cells (B4 date range, B5 download timestamp) ---
$ws.Range("B4").Value = "2024-11-01 ~ 2024-11-30"
$ws.Range("B5").Value = "2024년 12월 08일 16시 45분 54초"

# --- Drop the oldest daily-hit row (2024-10-01). Rows below shift up,
#     so the trailing blank spacer row lands back on the last row. ---
$ws.Rows(38).Delete()

# --- Refresh the daily-hit table (rows 8-37) for the new Nov-2024 window. ---
# Only cells whose value actually changes are rewritten; a leading "'" forces
# numeric-looking values (dates, counts) to stay stored as text, matching the
# original inlineStr cells instead of being reinterpreted as dates/numbers.
# row 8: 2024-10-31 -> 2024-11-30 (목 -> 토)
$ws.Range("A8").Value = "'2024-11-30"
$ws.Range("B8").Value = "토"
$ws.Range("C8").Value = "'9"
$ws.Range("F8").Value = "'9"

# row 9: 2024-10-30 -> 2024-11-29 (수 -> 금)
$ws.Range("A9").Value = "'2024-11-29"
$ws.Range("B9").Value = "금"
$ws.Range("C9").Value = "'8"
$ws.Range("D9").Value = "'0"
$ws.Range("F9").Value = "'8"

# row 10: 2024-10-29 -> 2024-11-28 (화 -> 목)
$ws.Range("A10").Value = "'2024-11-28"
$ws.Range("B10").Value = "목"
$ws.Range("C10").Value = "'8"
$ws.Range("F10").Value = "'8"

# row 11: 2024-10-28 -> 2024-11-27 (월 -> 수)
$ws.Range("A11").Value = "'2024-11-27"
$ws.Range("B11").Value = "수"
$ws.Range("C11").Value = "'10"
$ws.Range("F11").Value = "'10"

# row 12: 2024-10-27 -> 2024-11-26 (일 -> 화)
$ws.Range("A12").Value = "'2024-11-26"
$ws.Range("B12").Value = "화"
$ws.Range("C12").Value = "'20"
$ws.Range("F12").Value = "'20"

# row 13: 2024-10-26 -> 2024-11-25 (토 -> 월)
$ws.Range("A13").Value = "'2024-11-25"
$ws.Range("B13").Value = "월"

# row 14: 2024-10-25 -> 2024-11-24 (금 -> 일)
$ws.Range("A14").Value = "'2024-11-24"
$ws.Range("B14").Value = "일"
$ws.Range("C14").Value = "'7"
$ws.Range("F14").Value = "'7"

# row 15: 2024-10-24 -> 2024-11-23 (목 -> 토)
$ws.Range("A15").Value = "'2024-11-23"
$ws.Range("B15").Value = "토"
$ws.Range("C15").Value = "'12"
$ws.Range("E15").Value = "'1"
$ws.Range("F15").Value = "'11"

# row 16: 2024-10-23 -> 2024-11-22 (수 -> 금)
$ws.Range("A16").Value = "'2024-11-22"
$ws.Range("B16").Value = "금"
$ws.Range("C16").Value = "'7"
$ws.Range("F16").Value = "'7"

# row 17: 2024-10-22 -> 2024-11-21 (화 -> 목)
$ws.Range("A17").Value = "'2024-11-21"
$ws.Range("B17").Value = "목"
$ws.Range("C17").Value = "'9"
$ws.Range("F17").Value = "'9"

# row 18: 2024-10-21 -> 2024-11-20 (월 -> 수)
$ws.Range("A18").Value = "'2024-11-20"
$ws.Range("B18").Value = "수"
$ws.Range("C18").Value = "'8"
$ws.Range("D18").Value = "'1"
$ws.Range("F18").Value = "'7"

# row 19: 2024-10-20 -> 2024-11-19 (일 -> 화)
$ws.Range("A19").Value = "'2024-11-19"
$ws.Range("B19").Value = "화"
$ws.Range("C19").Value = "'8"
$ws.Range("F19").Value = "'8"

# row 20: 2024-10-19 -> 2024-11-18 (토 -> 월)
$ws.Range("A20").Value = "'2024-11-18"
$ws.Range("B20").Value = "월"
$ws.Range("C20").Value = "'7"
$ws.Range("D20").Value = "'2"
$ws.Range("F20").Value = "'5"

# row 21: 2024-10-18 -> 2024-11-17 (금 -> 일)
$ws.Range("A21").Value = "'2024-11-17"
$ws.Range("B21").Value = "일"

# row 22: 2024-10-17 -> 2024-11-16 (목 -> 토)
$ws.Range("A22").Value = "'2024-11-16"
$ws.Range("B22").Value = "토"
$ws.Range("C22").Value = "'10"
$ws.Range("F22").Value = "'10"

# row 23: 2024-10-16 -> 2024-11-15 (수 -> 금)
$ws.Range("A23").Value = "'2024-11-15"
$ws.Range("B23").Value = "금"
$ws.Range("C23").Value = "'16"
$ws.Range("D23").Value = "'5"
$ws.Range("F23").Value = "'11"

# row 24: 2024-10-15 -> 2024-11-14 (화 -> 목)
$ws.Range("A24").Value = "'2024-11-14"
$ws.Range("B24").Value = "목"
$ws.Range("C24").Value = "'13"
$ws.Range("D24").Value = "'8"
$ws.Range("F24").Value = "'5"

# row 25: 2024-10-14 -> 2024-11-13 (월 -> 수)
$ws.Range("A25").Value = "'2024-11-13"
$ws.Range("B25").Value = "수"
$ws.Range("C25").Value = "'21"
$ws.Range("D25").Value = "'1"
$ws.Range("F25").Value = "'20"

# row 26: 2024-10-13 -> 2024-11-12 (일 -> 화)
$ws.Range("A26").Value = "'2024-11-12"
$ws.Range("B26").Value = "화"
$ws.Range("C26").Value = "'7"
$ws.Range("D26").Value = "'1"
$ws.Range("F26").Value = "'6"

# row 27: 2024-10-12 -> 2024-11-11 (토 -> 월)
$ws.Range("A27").Value = "'2024-11-11"
$ws.Range("B27").Value = "월"
$ws.Range("C27").Value = "'21"
$ws.Range("D27").Value = "'0"
$ws.Range("E27").Value = "'0"
$ws.Range("F27").Value = "'21"

# row 28: 2024-10-11 -> 2024-11-10 (금 -> 일)
$ws.Range("A28").Value = "'2024-11-10"
$ws.Range("B28").Value = "일"
$ws.Range("C28").Value = "'12"
$ws.Range("D28").Value = "'1"
$ws.Range("E28").Value = "'1"
$ws.Range("F28").Value = "'10"

# row 29: 2024-10-10 -> 2024-11-09 (목 -> 토)
$ws.Range("A29").Value = "'2024-11-09"
$ws.Range("B29").Value = "토"
$ws.Range("C29").Value = "'13"
$ws.Range("D29").Value = "'0"
$ws.Range("E29").Value = "'1"
$ws.Range("F29").Value = "'12"

# row 30: 2024-10-09 -> 2024-11-08 (수 -> 금)
$ws.Range("A30").Value = "'2024-11-08"
$ws.Range("B30").Value = "금"
$ws.Range("C30").Value = "'17"
$ws.Range("D30").Value = "'4"
$ws.Range("E30").Value = "'0"
$ws.Range("F30").Value = "'13"

# row 31: 2024-10-08 -> 2024-11-07 (화 -> 목)
$ws.Range("A31").Value = "'2024-11-07"
$ws.Range("B31").Value = "목"
$ws.Range("C31").Value = "'12"
$ws.Range("D31").Value = "'3"
$ws.Range("F31").Value = "'9"

# row 32: 2024-10-07 -> 2024-11-06 (월 -> 수)
$ws.Range("A32").Value = "'2024-11-06"
$ws.Range("B32").Value = "수"
$ws.Range("C32").Value = "'12"
$ws.Range("F32").Value = "'12"

# row 33: 2024-10-06 -> 2024-11-05 (일 -> 화)
$ws.Range("A33").Value = "'2024-11-05"
$ws.Range("B33").Value = "화"
$ws.Range("C33").Value = "'8"
$ws.Range("F33").Value = "'8"

# row 34: 2024-10-05 -> 2024-11-04 (토 -> 월)
$ws.Range("A34").Value = "'2024-11-04"
$ws.Range("B34").Value = "월"
$ws.Range("C34").Value = "'9"
$ws.Range("F34").Value = "'9"

# row 35: 2024-10-04 -> 2024-11-03 (금 -> 일)
$ws.Range("A35").Value = "'2024-11-03"
$ws.Range("B35").Value = "일"
$ws.Range("C35").Value = "'7"
$ws.Range("F35").Value = "'7"

# row 36: 2024-10-03 -> 2024-11-02 (목 -> 토)
$ws.Range("A36").Value = "'2024-11-02"
$ws.Range("B36").Value = "토"
$ws.Range("D36").Value = "'0"
$ws.Range("F36").Value = "'13"

# row 37: 2024-10-02 -> 2024-11-01 (수 -> 금)
$ws.Range("A37").Value = "'2024-11-01"
$ws.Range("B37").Value = "금"
$ws.Range("C37").Value = "'4"
$ws.Range("F37").Value = "'4"

